$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Title text change (row 1)
$ws.Range("A1").Value = "Benchmark For Test/AND Vulnerability Factor"

# Column F header text change (row 2)
$ws.Range("F2").Value = "Gate Delay"

# Rows 3-8: "None" text values in column F become numeric 0
$ws.Range("F3:F8").Value = 0

# Rows 9-11: Gate inputs list text in column F becomes numeric gate delay 19.8
$ws.Range("F9").Value = 19.8
$ws.Range("F10").Value = 19.8
$ws.Range("F11").Value = 19.8

# Row 12: Gate inputs list text in column F becomes numeric gate delay 39.6
$ws.Range("F12").Value = 39.6
